$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 3-7 (columns D, I, J, K, L, M, N, P, Q) get cyclically rotated:
#   new row 3 <- old row 7
#   new row 4 <- old row 5
#   new row 5 <- old row 6
#   new row 6 <- old row 3
#   new row 7 <- old row 4
# Capture the "before" values for the affected columns first, then write them
# back out in rotated order so we don't clobber source data while writing.

$cols = @("D","I","J","K","L","M","N","P","Q")

$snapshot = @{}
foreach ($r in 3..7) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

$mapping = @{ 3 = 7; 4 = 5; 5 = 6; 6 = 3; 7 = 4 }

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
